$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "327.47"
Set-TextValue $ws.Range("E2") "-0.83%"
Set-TextValue $ws.Range("D3") "43.68"
Set-TextValue $ws.Range("E3") "5.05%"
Set-TextValue $ws.Range("D4") "5.490"
Set-TextValue $ws.Range("E4") "-3.36%"
Set-TextValue $ws.Range("D5") "0.08064"
Set-TextValue $ws.Range("E5") "-4.32%"
Set-TextValue $ws.Range("D6") "8.635"
Set-TextValue $ws.Range("E6") "-1.84%"
Set-TextValue $ws.Range("E7") "-4.31%"
Set-TextValue $ws.Range("D8") "1.874"
Set-TextValue $ws.Range("E8") "-5.73%"
Set-TextValue $ws.Range("E9") "-8.55%"
Set-TextValue $ws.Range("D10") "0.9374"
Set-TextValue $ws.Range("E10") "1.16%"
Set-TextValue $ws.Range("D11") "0.1189"
Set-TextValue $ws.Range("E11") "-6.88%"
Set-TextValue $ws.Range("D12") "0.1890"
Set-TextValue $ws.Range("E12") "-3.71%"
Set-TextValue $ws.Range("D13") "0.09565"
Set-TextValue $ws.Range("E13") "1.72%"
Set-TextValue $ws.Range("D14") "0.04121"
Set-TextValue $ws.Range("E14") "4.20%"
Set-TextValue $ws.Range("E15") "0.40%"
Set-TextValue $ws.Range("D16") "0.001281"
Set-TextValue $ws.Range("E16") "-2.54%"
Set-TextValue $ws.Range("D17") "0.005988"
Set-TextValue $ws.Range("E17") "-2.03%"
Set-TextValue $ws.Range("D18") "3.572"
Set-TextValue $ws.Range("E18") "4.34%"
Set-TextValue $ws.Range("D20") "8.557"
Set-TextValue $ws.Range("E20") "-4.24%"
Set-TextValue $ws.Range("D21") "0.1366"
Set-TextValue $ws.Range("E21") "0.14%"
Set-TextValue $ws.Range("D22") "0.2593"
Set-TextValue $ws.Range("E22") "3.18%"
Set-TextValue $ws.Range("D23") "0.04331"
Set-TextValue $ws.Range("E23") "-2.06%"
Set-TextValue $ws.Range("E24") "-0.98%"
Set-TextValue $ws.Range("D25") "0.004332"
Set-TextValue $ws.Range("E25") "-1.58%"
Set-TextValue $ws.Range("E26") "3.43%"
Set-TextValue $ws.Range("D27") "0.0003999"
Set-TextValue $ws.Range("E27") "0.12%"
Set-TextValue $ws.Range("D39") "0.02635"
Set-TextValue $ws.Range("E39") "-6.76%"
Set-TextValue $ws.Range("D40") "0.05400"
Set-TextValue $ws.Range("E40") "-2.04%"
Set-TextValue $ws.Range("D41") "0.01110"
Set-TextValue $ws.Range("E41") "23.49%"
Set-TextValue $ws.Range("D42") "0.007651"
Set-TextValue $ws.Range("E42") "-3.20%"
Set-TextValue $ws.Range("D43") "0.1388"
Set-TextValue $ws.Range("E43") "-3.41%"
Set-TextValue $ws.Range("E44") "2.11%"
Set-TextValue $ws.Range("D45") "0.009667"
Set-TextValue $ws.Range("E45") "-12.16%"
Set-TextValue $ws.Range("D46") "0.00006872"
Set-TextValue $ws.Range("E46") "-4.91%"
Set-TextValue $ws.Range("E47") "0.09%"
Set-TextValue $ws.Range("D48") "0.003569"
Set-TextValue $ws.Range("E48") "9.49%"
Set-TextValue $ws.Range("D49") "0.002276"
Set-TextValue $ws.Range("E49") "-0.22%"
Set-TextValue $ws.Range("D50") "0.00002105"
Set-TextValue $ws.Range("E50") "0.09%"
Set-TextValue $ws.Range("D51") "0.0002005"
Set-TextValue $ws.Range("E51") "0.09%"
